# Updates the auto-updating "datetimeFigureOut" date placeholder text
# (slide master + all slide layouts) from 15/12/2025 to 16/12/2025, and
# nudges the "GW-P 2020" rounded-rectangle shape on slide 1 to its new
# position.

$p = $ppt.ActivePresentation

$newDate = "16/12/2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    Update-DatePlaceholder $lay.Shapes
}

# Slide 1: move the "GW-P 2020" rounded-rectangle shape.
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(1)
$shp1.Left = 131.77371078740157
$shp1.Top = 137.2075590551181
